$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 92 is the last existing data row; append the new row 93 right after it
# with the same 4-column layout (date, weekday, hour, ranking).
#
# Column A holds a date-like string ("2025/10/11") that must stay plain text,
# matching every other date cell already in the sheet (they are all text,
# not real dates). A bare Value assignment would let Excel auto-convert the
# text to a date serial number, so a leading apostrophe forces text entry
# (classic "quote prefix" trick); ClearFormats() afterwards drops the
# quote-prefix formatting flag again so the cell is left with the sheet's
# default (unstyled) look, same as the surrounding data cells.
$ws.Range("A93").Value = "'2025/10/11"
$ws.Range("A93").ClearFormats()

$ws.Range("B93").Value = "土"
$ws.Range("C93").Value = 7
$ws.Range("D93").Value = 201
